$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell while forcing Excel to keep it as
# plain text (these columns hold formatted numeric/percentage strings like
# "243.04" or "-0.59%" that must NOT be auto-converted to numbers). Briefly
# switching to the "@" text format stops Excel re-interpreting the string as
# a number/percentage; restoring the "Normal" style afterwards keeps the
# cell's original (default) formatting/style untouched.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - BNB
Set-TextValue $ws.Range("D2") "243.04"
Set-TextValue $ws.Range("E2") "-0.59%"

# Row 3 - OKB
Set-TextValue $ws.Range("D3") "30.09"
Set-TextValue $ws.Range("E3") "13.74%"

# Row 4 - HuobiToken
Set-TextValue $ws.Range("D4") "5.142"
Set-TextValue $ws.Range("E4") "0.30%"

# Row 5 - Cronos
Set-TextValue $ws.Range("D5") "0.05672"
Set-TextValue $ws.Range("E5") "1.52%"

# Row 6
Set-TextValue $ws.Range("D6") "6.518"
Set-TextValue $ws.Range("E6") "0.79%"

# Row 7
Set-TextValue $ws.Range("D7") "0.8423"
Set-TextValue $ws.Range("E7") "3.07%"

# Row 8
Set-TextValue $ws.Range("D8") "0.8630"
Set-TextValue $ws.Range("E8") "3.39%"

# Row 9
Set-TextValue $ws.Range("D9") "0.1338"
Set-TextValue $ws.Range("E9") "0.39%"

# Row 10
Set-TextValue $ws.Range("D10") "0.06917"
Set-TextValue $ws.Range("E10") "-1.14%"

# Row 11
Set-TextValue $ws.Range("D11") "0.02858"
Set-TextValue $ws.Range("E11") "-1.03%"

# Row 12
Set-TextValue $ws.Range("D12") "0.09381"
Set-TextValue $ws.Range("E12") "-0.06%"

# Row 13
Set-TextValue $ws.Range("D13") "0.001517"
Set-TextValue $ws.Range("E13") "0.63%"

# Row 14 (only E changes)
Set-TextValue $ws.Range("E14") "-11.01%"

# Row 15
Set-TextValue $ws.Range("D15") "0.0005975"
Set-TextValue $ws.Range("E15") "0.41%"

# Row 16
Set-TextValue $ws.Range("D16") "0.006081"
Set-TextValue $ws.Range("E16") "-1.13%"

# Row 17
Set-TextValue $ws.Range("D17") "3.509"
Set-TextValue $ws.Range("E17") "-3.81%"

# Row 18
Set-TextValue $ws.Range("D18") "3.020"
Set-TextValue $ws.Range("E18") "-0.53%"

# Row 19 (only E changes)
Set-TextValue $ws.Range("E19") "-2.55%"

# Row 20
Set-TextValue $ws.Range("D20") "0.3152"
Set-TextValue $ws.Range("E20") "1.28%"

# Row 21
Set-TextValue $ws.Range("D21") "0.03267"
Set-TextValue $ws.Range("E21") "5.21%"

# Row 22 (only E changes)
Set-TextValue $ws.Range("E22") "-0.32%"

# Row 23
Set-TextValue $ws.Range("D23") "3.615"
Set-TextValue $ws.Range("E23") "-3.31%"

# Row 24 (only E changes)
Set-TextValue $ws.Range("E24") "-0.11%"

# Row 25
Set-TextValue $ws.Range("D25") "0.001210"
Set-TextValue $ws.Range("E25") "-2.72%"

# Row 26
Set-TextValue $ws.Range("D26") "0.004321"
Set-TextValue $ws.Range("E26") "-4.00%"

# Row 27
Set-TextValue $ws.Range("D27") "0.0001180"
Set-TextValue $ws.Range("E27") "22.89%"

# Row 28 (only E changes)
Set-TextValue $ws.Range("E28") "0.23%"

# Row 40
Set-TextValue $ws.Range("D40") "0.03712"
Set-TextValue $ws.Range("E40") "1.93%"

# Row 41 - becomes BKEXToken (was KickToken)
Set-TextValue $ws.Range("B41") "BKEXToken"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D41") "0.1058"
Set-TextValue $ws.Range("E41") "0.79%"

# Row 42 - becomes KickToken (was BKEXToken)
Set-TextValue $ws.Range("B42") "KickToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D42") "0.005329"
Set-TextValue $ws.Range("E42") "-13.63%"

# Row 43
Set-TextValue $ws.Range("D43") "0.002310"
Set-TextValue $ws.Range("E43") "-3.75%"

# Row 44
Set-TextValue $ws.Range("D44") "0.009737"
Set-TextValue $ws.Range("E44") "9.76%"

# Row 45
Set-TextValue $ws.Range("D45") "0.00005093"
Set-TextValue $ws.Range("E45") "-4.99%"

# Row 46 (only E changes)
Set-TextValue $ws.Range("E46") "-0.09%"

# Row 47
Set-TextValue $ws.Range("D47") "0.09991"
Set-TextValue $ws.Range("E47") "-30.62%"

# Row 48
Set-TextValue $ws.Range("D48") "0.002720"
Set-TextValue $ws.Range("E48") "16.47%"

# Row 49 (only E changes)
Set-TextValue $ws.Range("E49") "-0.09%"

# Row 50 (only E changes)
Set-TextValue $ws.Range("E50") "-0.09%"
